# "added laptop benchmarks and serial flag -s"
# Re-times the second ("Laptop") benchmark table, reformats the timing
# columns to 6 decimal places, inserts a spacer row above the second
# table's header, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the new 6-decimal number format to the timing columns -----
# (first table E4:F9, the new spacer row E10:F10, the second table's
# header row E11:F11, and the second table's data E12:F17)
$ws.Range("E4:F9").NumberFormat = "0.000000"
$ws.Range("E10:F10").NumberFormat = "0.000000"
$ws.Range("E11:F11").NumberFormat = "0.000000"
$ws.Range("E12:F17").NumberFormat = "0.000000"

# --- New laptop timings (re-measured serial run, flag -s) ------------
$ws.Range("E12").Value = 0.349148
$ws.Range("F12").Value = 0.087414
$ws.Range("G12").Formula = "=E12/F12"

$ws.Range("E13").Value = 1.969367
$ws.Range("F13").Value = 0.249261
# This value got pasted over the old shared formula by hand, so G13 is a
# plain literal now instead of "=E13/F13".
$ws.Range("G13").Value = 0.359487

$ws.Range("E14").Value = 7.500394
$ws.Range("F14").Value = 1.363164
$ws.Range("G14").Formula = "=E14/F14"

$ws.Range("E15").Value = 1.043635
$ws.Range("F15").Value = 0.201269
$ws.Range("G15").Formula = "=E15/F15"

$ws.Range("E16").Value = 5.969409
$ws.Range("F16").Value = 0.983853
$ws.Range("G16").Formula = "=E16/F16"

$ws.Range("E17").Value = 22.80116
$ws.Range("F17").Value = 3.634914
$ws.Range("G17").Formula = "=E17/F17"

# --- Move the active selection ----------------------------------------
$ws.Range("H15").Select() | Out-Null
